# Applies the updated betting-odds values from the 2024-12-13 FlashScore refresh.
# Only the specific cells listed below change; everything else is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7
$ws.Cells.Item(7, 7).Value = 6  # G7: 6.25 -> 6
$ws.Cells.Item(7, 8).Value = 3.75  # H7: 3.7 -> 3.75
$ws.Cells.Item(7, 13).Value = 1.08  # M7: 1.1 -> 1.08
$ws.Cells.Item(7, 14).Value = 8  # N7: 7 -> 8
$ws.Cells.Item(7, 15).Value = 1.44  # O7: 1.5 -> 1.44
$ws.Cells.Item(7, 16).Value = 2.63  # P7: 2.5 -> 2.63
$ws.Cells.Item(7, 17).Value = 2.4  # Q7: 2.5 -> 2.4
$ws.Cells.Item(7, 18).Value = 1.53  # R7: 1.5 -> 1.53
$ws.Cells.Item(7, 19).Value = 1.53  # S7: 1.57 -> 1.53
$ws.Cells.Item(7, 20).Value = 2.38  # T7: 2.25 -> 2.38
$ws.Cells.Item(7, 26).Value = 67  # Z7: 81 -> 67
$ws.Cells.Item(7, 29).Value = 7  # AC7: 6.5 -> 7
$ws.Cells.Item(7, 43).Value = 151  # AQ7: 201 -> 151
$ws.Cells.Item(7, 46).Value = 2.38  # AT7: 2.25 -> 2.38
$ws.Cells.Item(7, 47).Value = 10  # AU7: 11 -> 10
$ws.Cells.Item(7, 54).Value = 251  # BB7: 301 -> 251
$ws.Cells.Item(7, 56).Value = 151  # BD7: 126 -> 151

# Row 8
$ws.Cells.Item(8, 7).Value = 2.4  # G8: 2.45 -> 2.4
$ws.Cells.Item(8, 9).Value = 3.2  # I8: 3.1 -> 3.2
$ws.Cells.Item(8, 10).Value = 3.2  # J8: 3.25 -> 3.2
$ws.Cells.Item(8, 12).Value = 4  # L8: 3.75 -> 4
$ws.Cells.Item(8, 24).Value = 10  # X8: 11 -> 10
$ws.Cells.Item(8, 35).Value = 15  # AI8: 13 -> 15

# Row 9
$ws.Cells.Item(9, 7).Value = 2.3  # G9: 2.4 -> 2.3
$ws.Cells.Item(9, 9).Value = 3.5  # I9: 3.3 -> 3.5
$ws.Cells.Item(9, 10).Value = 3.2  # J9: 3.4 -> 3.2
$ws.Cells.Item(9, 12).Value = 4.5  # L9: 4.33 -> 4.5
$ws.Cells.Item(9, 13).Value = 1.14  # M9: 1.13 -> 1.14
$ws.Cells.Item(9, 14).Value = 5.5  # N9: 6 -> 5.5
$ws.Cells.Item(9, 24).Value = 9  # X9: 9.5 -> 9
$ws.Cells.Item(9, 26).Value = 21  # Z9: 23 -> 21
$ws.Cells.Item(9, 34).Value = 7  # AH9: 6.5 -> 7
$ws.Cells.Item(9, 36).Value = 15  # AJ9: 13 -> 15
$ws.Cells.Item(9, 38).Value = 41  # AL9: 34 -> 41
$ws.Cells.Item(9, 40).Value = 4  # AN9: 4.33 -> 4
$ws.Cells.Item(9, 49).Value = 5.5  # AW9: 5 -> 5.5
$ws.Cells.Item(9, 50).Value = 23  # AX9: 21 -> 23
$ws.Cells.Item(9, 54).Value = 501  # BB9: 451 -> 501

# Row 10
$ws.Cells.Item(10, 9).Value = 4.33  # I10: 4.5 -> 4.33
$ws.Cells.Item(10, 10).Value = 2.75  # J10: 2.88 -> 2.75
$ws.Cells.Item(10, 11).Value = 1.95  # K10: 1.91 -> 1.95
$ws.Cells.Item(10, 17).Value = 2.5  # Q10: 2.6 -> 2.5
$ws.Cells.Item(10, 18).Value = 1.5  # R10: 1.48 -> 1.5
$ws.Cells.Item(10, 27).Value = 19  # AA10: 21 -> 19
$ws.Cells.Item(10, 29).Value = 6.5  # AC10: 6 -> 6.5
$ws.Cells.Item(10, 49).Value = 5.5  # AW10: 6 -> 5.5

# Row 15
$ws.Cells.Item(15, 7).Value = 2.9  # G15: 3 -> 2.9
$ws.Cells.Item(15, 9).Value = 2.35  # I15: 2.3 -> 2.35
$ws.Cells.Item(15, 10).Value = 3.4  # J15: 3.5 -> 3.4
$ws.Cells.Item(15, 12).Value = 3  # L15: 2.88 -> 3
$ws.Cells.Item(15, 17).Value = 1.98  # Q15: 1.95 -> 1.98
$ws.Cells.Item(15, 18).Value = 1.83  # R15: 1.85 -> 1.83
$ws.Cells.Item(15, 26).Value = 29  # Z15: 34 -> 29
$ws.Cells.Item(15, 28).Value = 29  # AB15: 34 -> 29
$ws.Cells.Item(15, 29).Value = 11  # AC15: 10 -> 11
$ws.Cells.Item(15, 34).Value = 8.5  # AH15: 8 -> 8.5
$ws.Cells.Item(15, 35).Value = 12  # AI15: 11 -> 12
$ws.Cells.Item(15, 37).Value = 23  # AK15: 21 -> 23
$ws.Cells.Item(15, 41).Value = 15  # AO15: 17 -> 15
$ws.Cells.Item(15, 50).Value = 13  # AX15: 12 -> 13

# Row 16
$ws.Cells.Item(16, 7).Value = 1.14  # G16: 1.11 -> 1.14
$ws.Cells.Item(16, 8).Value = 8.5  # H16: 9.5 -> 8.5
$ws.Cells.Item(16, 9).Value = 13  # I16: 15 -> 13
$ws.Cells.Item(16, 10).Value = 1.5  # J16: 1.4 -> 1.5
$ws.Cells.Item(16, 11).Value = 3.25  # K16: 3.5 -> 3.25
$ws.Cells.Item(16, 12).Value = 10  # L16: 12 -> 10
$ws.Cells.Item(16, 15).Value = 1.07  # O16: 1.06 -> 1.07
$ws.Cells.Item(16, 16).Value = 9  # P16: 10 -> 9
$ws.Cells.Item(16, 17).Value = 1.29  # Q16: 1.25 -> 1.29
$ws.Cells.Item(16, 18).Value = 3.6  # R16: 4 -> 3.6
$ws.Cells.Item(16, 21).Value = 1.83  # U16: 1.91 -> 1.83
$ws.Cells.Item(16, 22).Value = 1.83  # V16: 1.8 -> 1.83
$ws.Cells.Item(16, 25).Value = 11  # Y16: 12 -> 11
$ws.Cells.Item(16, 26).Value = 8  # Z16: 7.5 -> 8
$ws.Cells.Item(16, 27).Value = 10  # AA16: 11 -> 10
$ws.Cells.Item(16, 28).Value = 23  # AB16: 26 -> 23
$ws.Cells.Item(16, 29).Value = 29  # AC16: 34 -> 29
$ws.Cells.Item(16, 30).Value = 17  # AD16: 19 -> 17
$ws.Cells.Item(16, 31).Value = 26  # AE16: 29 -> 26
$ws.Cells.Item(16, 32).Value = 51  # AF16: 67 -> 51
$ws.Cells.Item(16, 35).Value = 67  # AI16: 81 -> 67
$ws.Cells.Item(16, 36).Value = 34  # AJ16: 41 -> 34
$ws.Cells.Item(16, 37).Value = 151  # AK16: 201 -> 151
$ws.Cells.Item(16, 38).Value = 67  # AL16: 81 -> 67
$ws.Cells.Item(16, 39).Value = 51  # AM16: 67 -> 51
$ws.Cells.Item(16, 40).Value = 3.5  # AN16: 3.6 -> 3.5
$ws.Cells.Item(16, 41).Value = 5  # AO16: 4.75 -> 5
$ws.Cells.Item(16, 42).Value = 13  # AP16: 15 -> 13
$ws.Cells.Item(16, 43).Value = 10  # AQ16: 9 -> 10
$ws.Cells.Item(16, 47).Value = 10  # AU16: 11 -> 10
$ws.Cells.Item(16, 49).Value = 13  # AW16: 15 -> 13
$ws.Cells.Item(16, 52).Value = 251  # AZ16: 301 -> 251
$ws.Cells.Item(16, 54).Value = 251  # BB16: 301 -> 251

# Row 17
$ws.Cells.Item(17, 7).Value = 1.27  # G17: 1.25 -> 1.27
$ws.Cells.Item(17, 9).Value = 10  # I17: 11 -> 10
$ws.Cells.Item(17, 15).Value = 1.11  # O17: 1.13 -> 1.11
$ws.Cells.Item(17, 16).Value = 6.5  # P17: 6 -> 6.5
$ws.Cells.Item(17, 17).Value = 1.4  # Q17: 1.44 -> 1.4
$ws.Cells.Item(17, 18).Value = 2.88  # R17: 2.7 -> 2.88
$ws.Cells.Item(17, 23).Value = 9.5  # W17: 10 -> 9.5
$ws.Cells.Item(17, 24).Value = 8  # X17: 7.5 -> 8
$ws.Cells.Item(17, 25).Value = 9  # Y17: 9.5 -> 9
$ws.Cells.Item(17, 28).Value = 21  # AB17: 23 -> 21
$ws.Cells.Item(17, 30).Value = 11  # AD17: 12 -> 11
$ws.Cells.Item(17, 31).Value = 19  # AE17: 21 -> 19
$ws.Cells.Item(17, 34).Value = 34  # AH17: 29 -> 34
$ws.Cells.Item(17, 41).Value = 6  # AO17: 5.5 -> 6
$ws.Cells.Item(17, 45).Value = 81  # AS17: 101 -> 81
$ws.Cells.Item(17, 51).Value = 34  # AY17: 41 -> 34
$ws.Cells.Item(17, 53).Value = 126  # BA17: 151 -> 126
$ws.Cells.Item(17, 54).Value = 201  # BB17: 251 -> 201
